# Updates the HLN/DM statistical-ranking workbook ("cambios en la adicion"):
# the underlying HLN/DM comparison was recomputed, which (a) re-ranks a few
# models on the "Ranking" sheet (win/tie/score counts + CRPS stats, and the
# model names attached to a handful of rows shift accordingly) and
# (b) refreshes the superiority / p-value / DM-statistic matrices that back
# it. Apply the new cell values directly.
$wb = $excel.ActiveWorkbook

# --- Sheet "Ranking": re-sorted standings, updated win/tie/score + CRPS ---
$ws = $wb.Worksheets.Item("Ranking")

# Row 2 (Sieve Bootstrap): ties/score/win-rate recomputed
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.125

# Row 4 now holds EnCQR-LSTM (was AV-MCPS)
$ws.Range("B4").Value = "EnCQR-LSTM"
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.125
$ws.Range("H4").Value = 0.002862421887232199
$ws.Range("I4").Value = 0.002737018741784631

# Row 5 now holds AV-MCPS (was LSPMW)
$ws.Range("B5").Value = "AV-MCPS"
$ws.Range("H5").Value = 0.00233372413512631
$ws.Range("I5").Value = 0.001695838571145866

# Row 7 (MCPS): CRPS stats refreshed
$ws.Range("H7").Value = 0.002647994599202048
$ws.Range("I7").Value = 0.002238418271929144

# Row 8 now holds LSPMW (was Block Bootstrapping)
$ws.Range("B8").Value = "LSPMW"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 0
$ws.Range("H8").Value = 0.01088643503280608
$ws.Range("I8").Value = 0.01118631828981691

# Row 9 now holds Block Bootstrapping (was AREPD)
$ws.Range("B9").Value = "Block Bootstrapping"
$ws.Range("H9").Value = 0.01074678080948211
$ws.Range("I9").Value = 0.01054280575249088

# Row 10 now holds AREPD (was EnCQR-LSTM)
$ws.Range("B10").Value = "AREPD"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = -2
$ws.Range("H10").Value = 0.01085869601634492
$ws.Range("I10").Value = 0.01107257398892365

# --- Sheet "Matriz_Superioridad": superiority flags vs. EnCQR-LSTM updated ---
$ws = $wb.Worksheets.Item("Matriz_Superioridad")
$ws.Range("J3").Value = 0
$ws.Range("J6").Value = -1
$ws.Range("C10").Value = 0
$ws.Range("F10").Value = 1

# --- Sheet "Matriz_Pvalores": p-values recomputed (symmetric matrix) ---
$ws = $wb.Worksheets.Item("Matriz_Pvalores")
$ws.Range("G2").Value = 0.003174807860589546
$ws.Range("H2").Value = 0.004079998624447523
$ws.Range("J2").Value = 0.004439147514800013
$ws.Range("G3").Value = 0.2527363784636867
$ws.Range("H3").Value = 0.3711699617564306
$ws.Range("J3").Value = 0.02644120338402178
$ws.Range("G4").Value = 0.1154261312600404
$ws.Range("H4").Value = 0.08419339035355633
$ws.Range("J4").Value = 0.1329130947100663
$ws.Range("G5").Value = 0.02716039587005836
$ws.Range("H5").Value = 0.02884605071504298
$ws.Range("J5").Value = 0.02938501293574203
$ws.Range("G6").Value = 0.002377401658132605
$ws.Range("H6").Value = 0.002445029713906255
$ws.Range("J6").Value = 0.001171142210408815
$ws.Range("B7").Value = 0.003174807860589546
$ws.Range("C7").Value = 0.2527363784636867
$ws.Range("D7").Value = 0.1154261312600404
$ws.Range("E7").Value = 0.02716039587005836
$ws.Range("F7").Value = 0.002377401658132605
$ws.Range("H7").Value = 0.388923800022529
$ws.Range("I7").Value = 0.5237353623974466
$ws.Range("J7").Value = 0.7323148542116584
$ws.Range("B8").Value = 0.004079998624447523
$ws.Range("C8").Value = 0.3711699617564306
$ws.Range("D8").Value = 0.08419339035355633
$ws.Range("E8").Value = 0.02884605071504298
$ws.Range("F8").Value = 0.002445029713906255
$ws.Range("G8").Value = 0.388923800022529
$ws.Range("I8").Value = 0.381143703974532
$ws.Range("J8").Value = 0.2927194695972357
$ws.Range("G9").Value = 0.5237353623974466
$ws.Range("H9").Value = 0.381143703974532
$ws.Range("J9").Value = 0.7566060207781393
$ws.Range("B10").Value = 0.004439147514800013
$ws.Range("C10").Value = 0.02644120338402178
$ws.Range("D10").Value = 0.1329130947100663
$ws.Range("E10").Value = 0.02938501293574203
$ws.Range("F10").Value = 0.001171142210408815
$ws.Range("G10").Value = 0.7323148542116584
$ws.Range("H10").Value = 0.2927194695972357
$ws.Range("I10").Value = 0.7566060207781393

# --- Sheet "Matriz_DM_Original": raw DM statistics recomputed (antisymmetric matrix) ---
$ws = $wb.Worksheets.Item("Matriz_DM_Original")
$ws.Range("G2").Value = 6.337129806051893
$ws.Range("H2").Value = 5.919247596615206
$ws.Range("J2").Value = 5.783987257324997
$ws.Range("G3").Value = -1.335180262293603
$ws.Range("H3").Value = -1.006381502477831
$ws.Range("J3").Value = -3.433891055953135
$ws.Range("G4").Value = 2.005224285639492
$ws.Range("H4").Value = 2.286367445133264
$ws.Range("J4").Value = 1.882393605666372
$ws.Range("G5").Value = 3.40466490293299
$ws.Range("H5").Value = 3.339617146343697
$ws.Range("J5").Value = 3.319763837016656
$ws.Range("G6").Value = 6.849876913593199
$ws.Range("H6").Value = 6.798634165467775
$ws.Range("J6").Value = 8.26166467492926
$ws.Range("B7").Value = -6.337129806051893
$ws.Range("C7").Value = 1.335180262293603
$ws.Range("D7").Value = -2.005224285639492
$ws.Range("E7").Value = -3.40466490293299
$ws.Range("F7").Value = -6.849876913593199
$ws.Range("H7").Value = 0.9656024611261345
$ws.Range("I7").Value = -0.6978033078871785
$ws.Range("J7").Value = -0.366848143389047
$ws.Range("B8").Value = -5.919247596615206
$ws.Range("C8").Value = 1.006381502477831
$ws.Range("D8").Value = -2.286367445133264
$ws.Range("E8").Value = -3.339617146343697
$ws.Range("F8").Value = -6.798634165467775
$ws.Range("G8").Value = -0.9656024611261345
$ws.Range("I8").Value = -0.9832701938040602
$ws.Range("J8").Value = -1.210506061369819
$ws.Range("G9").Value = 0.6978033078871785
$ws.Range("H9").Value = 0.9832701938040602
$ws.Range("J9").Value = 0.3319302933904444
$ws.Range("B10").Value = -5.783987257324997
$ws.Range("C10").Value = 3.433891055953135
$ws.Range("D10").Value = -1.882393605666372
$ws.Range("E10").Value = -3.319763837016656
$ws.Range("F10").Value = -8.26166467492926
$ws.Range("G10").Value = 0.366848143389047
$ws.Range("H10").Value = 1.210506061369819
$ws.Range("I10").Value = -0.3319302933904444
